$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1964.2
$ws.Range("J17").Value = 2264.5
$ws.Range("L17").Value = 6793.5
$ws.Range("N17").Value = -7129.5
$ws.Range("H116").Value = 4967.143
$ws.Range("I116").Value = 4720.273
$ws.Range("J116").Value = 5238.7
$ws.Range("K116").Value = 4720.273
$ws.Range("L116").Value = 5238.7
$ws.Range("M116").Value = -1278.273
$ws.Range("N116").Value = -12122.7
$ws.Range("H137").Value = 1595.8125
$ws.Range("J137").Value = 2201.4
$ws.Range("L137").Value = 6604.200000000001
$ws.Range("N137").Value = -11704.2
$ws.Range("H141").Value = 15864.583
$ws.Range("I141").Value = 5430.6665
$ws.Range("K141").Value = 16291.9995
$ws.Range("M141").Value = -11111.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 13872
$ws.Range("I16").Value = 13872
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 13872
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -13585
$ws.Range("N16").ClearContents()
$ws.Range("H26").Value = 12660.667
$ws.Range("I26").Value = 12660.667
$ws.Range("K26").Value = 12660.667
$ws.Range("M26").Value = -12330.667
$ws.Range("H32").Value = 4090.3618
$ws.Range("I32").Value = 2713.805
$ws.Range("K32").Value = 2713.805
$ws.Range("M32").Value = -2426.805
$ws.Range("H61").Value = 3770.75
$ws.Range("I61").Value = 2971.9092
$ws.Range("J61").Value = 5528.2
$ws.Range("K61").Value = 2971.9092
$ws.Range("L61").Value = 5528.2
$ws.Range("M61").Value = -2759.9092
$ws.Range("N61").Value = -5952.2
$ws.Range("H74").Value = 5154.8
$ws.Range("I74").Value = 5044
$ws.Range("K74").Value = 5044
$ws.Range("M74").Value = -4170
$ws.Range("H76").Value = 147296
$ws.Range("J76").Value = 147296
$ws.Range("L76").Value = 147296
$ws.Range("N76").Value = -147972
$ws.Range("H77").Value = 5154.8
$ws.Range("I77").Value = 5044
$ws.Range("K77").Value = 25220
$ws.Range("M77").Value = -20852
$ws.Range("H79").Value = 147296
$ws.Range("J79").Value = 147296
$ws.Range("L79").Value = 147296
$ws.Range("N79").Value = -149636
$ws.Range("H132").Value = 34338.832
$ws.Range("I132").Value = 10253.375
$ws.Range("K132").Value = 30760.125
$ws.Range("M132").Value = -28230.125
$ws.Range("H136").Value = 3770.75
$ws.Range("I136").Value = 2971.9092
$ws.Range("J136").Value = 5528.2
$ws.Range("K136").Value = 8915.7276
$ws.Range("L136").Value = 16584.6
$ws.Range("M136").Value = -6365.7276
$ws.Range("N136").Value = -21684.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 57444.625
$ws.Range("J92").Value = 57444.625
$ws.Range("L92").Value = 57444.625
$ws.Range("N92").Value = -62436.625
$ws.Range("H134").Value = 853.4706
$ws.Range("I134").Value = 844.3125
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 2532.9375
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = 2.0625
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 112.44444
$ws.Range("I7").Value = 103.36364
$ws.Range("K7").Value = 103.36364
$ws.Range("M7").Value = 9.636359999999996
$ws.Range("H19").Value = 1674.5454
$ws.Range("J19").Value = 950
$ws.Range("L19").Value = 950
$ws.Range("N19").Value = -1290
$ws.Range("H24").Value = 1674.5454
$ws.Range("J24").Value = 950
$ws.Range("L24").Value = 950
$ws.Range("N24").Value = -1290
$ws.Range("H31").Value = 4469.5
$ws.Range("J31").Value = 5838.2
$ws.Range("L31").Value = 5838.2
$ws.Range("N31").Value = -6428.2
$ws.Range("H34").Value = 4469.5
$ws.Range("J34").Value = 5838.2
$ws.Range("L34").Value = 5838.2
$ws.Range("N34").Value = -6242.2
$ws.Range("H74").Value = 74981.336
$ws.Range("J74").Value = 74981.336
$ws.Range("L74").Value = 74981.336
$ws.Range("N74").Value = -76729.336
$ws.Range("H77").Value = 74981.336
$ws.Range("J77").Value = 74981.336
$ws.Range("L77").Value = 224944.008
$ws.Range("N77").Value = -233680.008
$ws.Range("H94").Value = 1517.4
$ws.Range("J94").Value = 1019
$ws.Range("L94").Value = 1019
$ws.Range("N94").Value = -1921
$ws.Range("H99").Value = 1999.625
$ws.Range("I99").Value = 1999
$ws.Range("J99").Value = 1999.7142
$ws.Range("K99").Value = 1999
$ws.Range("L99").Value = 1999.7142
$ws.Range("M99").Value = -501
$ws.Range("N99").Value = -4995.7142
$ws.Range("H126").Value = 1999.625
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 1999.7142
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 5999.142599999999
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -10939.1426
$ws.Range("H132").Value = 2693.1667
$ws.Range("I132").Value = 2413.3333
$ws.Range("K132").Value = 7239.999899999999
$ws.Range("M132").Value = -4709.999899999999
$ws.Range("H134").Value = 2885.1538
$ws.Range("I134").Value = 2975.5833
$ws.Range("K134").Value = 8926.749899999999
$ws.Range("M134").Value = -6391.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 92381.45
$ws.Range("I4").Value = 100619.6
$ws.Range("K4").Value = 301858.8
$ws.Range("M4").Value = -301746.8
$ws.Range("H8").Value = 2324.6667
$ws.Range("I8").Value = 2324.6667
$ws.Range("K8").Value = 6974.000100000001
$ws.Range("M8").Value = -6835.000100000001
$ws.Range("H19").Value = 666
$ws.Range("J19").Value = 999
$ws.Range("L19").Value = 2997
$ws.Range("N19").Value = -3345
$ws.Range("H131").Value = 1659.5294
$ws.Range("J131").Value = 3106.6667
$ws.Range("L131").Value = 9320.000100000001
$ws.Range("N131").Value = -19400.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 29500.334
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1443.4615
$ws.Range("I46").Value = 1127.5
$ws.Range("K46").Value = 1127.5
$ws.Range("M46").Value = -939.5
$ws.Range("H64").Value = 79000
$ws.Range("J64").Value = 79000
$ws.Range("L64").Value = 79000
$ws.Range("N64").Value = -79450
$ws.Range("H67").Value = 79000
$ws.Range("J67").Value = 79000
$ws.Range("L67").Value = 79000
$ws.Range("N67").Value = -80560
$ws.Range("H122").Value = 11067.5
$ws.Range("I122").Value = 10997.272
$ws.Range("J122").Value = 11325
$ws.Range("K122").Value = 32991.81600000001
$ws.Range("L122").Value = 33975
$ws.Range("M122").Value = -30541.81600000001
$ws.Range("N122").Value = -38875
$ws.Range("H132").Value = 5120.644
$ws.Range("I132").Value = 4616.227
$ws.Range("K132").Value = 13848.681
$ws.Range("M132").Value = -11318.681
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059
$ws.Range("H136").Value = 6413.436
$ws.Range("I136").Value = 6003.3784
$ws.Range("K136").Value = 18010.1352
$ws.Range("M136").Value = -15460.1352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 79900
$ws.Range("J80").Value = 79900
$ws.Range("L80").Value = 79900
$ws.Range("N80").Value = -81896
$ws.Range("H83").Value = 79900
$ws.Range("J83").Value = 79900
$ws.Range("L83").Value = 239700
$ws.Range("N83").Value = -249684
$ws.Range("H97").Value = 52500
$ws.Range("J97").Value = 52500
$ws.Range("L97").Value = 52500
$ws.Range("N97").Value = -54482
$ws.Range("H109").Value = 60377
$ws.Range("J109").Value = 60377
$ws.Range("L109").Value = 60377
$ws.Range("N109").Value = -63151
$ws.Range("H122").Value = 3606.652
$ws.Range("I122").Value = 3148.5625
$ws.Range("J122").Value = 4653.7144
$ws.Range("K122").Value = 9445.6875
$ws.Range("L122").Value = 13961.1432
$ws.Range("M122").Value = -6995.6875
$ws.Range("N122").Value = -18861.1432
$ws.Range("H126").Value = 2640
$ws.Range("I126").Value = 2640
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7920
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5450
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1969.8334
$ws.Range("I132").Value = 1969.8334
$ws.Range("K132").Value = 5909.5002
$ws.Range("M132").Value = -3379.5002
$ws.Range("H136").Value = 3160.375
$ws.Range("I136").Value = 3210.9614
$ws.Range("K136").Value = 9632.8842
$ws.Range("L136").Value = 8823.500100000001
$ws.Range("M136").Value = -7082.8842
$ws.Range("N136").Value = -13923.5001
